# Update "想去人数" (want-to-go count) values in column F across sheets
# 展览 (sheet "展览"), 演出 (sheet "演出"), and 全部类型 (sheet "全部类型")

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 518
$ws1.Range("F3").Value = 739
$ws1.Range("F4").Value = 1465
$ws1.Range("F6").Value = 91
$ws1.Range("F7").Value = 138
$ws1.Range("F8").Value = 6171
$ws1.Range("F9").Value = 69
$ws1.Range("F10").Value = 401
$ws1.Range("F11").Value = 112
$ws1.Range("F12").Value = 5046
$ws1.Range("F13").Value = 25
$ws1.Range("F15").Value = 1172
$ws1.Range("F16").Value = 53
$ws1.Range("F18").Value = 61
$ws1.Range("F20").Value = 289
$ws1.Range("F21").Value = 25
$ws1.Range("F22").Value = 3569
$ws1.Range("F23").Value = 146

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 73

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 73
$ws4.Range("F3").Value = 518
$ws4.Range("F4").Value = 739
$ws4.Range("F5").Value = 1465
$ws4.Range("F7").Value = 91
$ws4.Range("F8").Value = 138
$ws4.Range("F9").Value = 6171
$ws4.Range("F10").Value = 69
$ws4.Range("F11").Value = 401
$ws4.Range("F12").Value = 112
$ws4.Range("F13").Value = 5046
$ws4.Range("F14").Value = 25
$ws4.Range("F16").Value = 1172
$ws4.Range("F17").Value = 53
$ws4.Range("F19").Value = 61
$ws4.Range("F21").Value = 289
$ws4.Range("F22").Value = 25
$ws4.Range("F23").Value = 3569
$ws4.Range("F25").Value = 146
